$wb = $excel.ActiveWorkbook
$hospital = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Data validation on the "Hospital" sheet (EDAD whole 0-125, VALOR CONSULTA
#    decimal >= 0)
# ---------------------------------------------------------------------------
$hospital.Range("D1:D1048576").Validation.Add(1, 1, 1, 0, 125) | Out-Null
$hospital.Range("E1:E1048576").Validation.Add(2, 1, 7, 0) | Out-Null

# Column widths on "Hospital" (best-fit to content, as after typing the data)
$hospital.Columns.Item(1).AutoFit() | Out-Null
$hospital.Columns.Item(2).AutoFit() | Out-Null
$hospital.Columns.Item(3).AutoFit() | Out-Null
$hospital.Columns.Item(4).AutoFit() | Out-Null
$hospital.Columns.Item(5).AutoFit() | Out-Null
$hospital.Columns.Item(6).AutoFit() | Out-Null
$hospital.Columns.Item(7).AutoFit() | Out-Null
$hospital.Columns.Item(8).AutoFit() | Out-Null
$hospital.Columns.Item(9).AutoFit() | Out-Null
$hospital.Columns.Item(10).AutoFit() | Out-Null
$hospital.Columns.Item(11).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 2. New worksheet "estadisticas" placed right after "Hospital", which also
#    becomes the active sheet/tab
# ---------------------------------------------------------------------------
$stats = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $hospital)
$stats.Name = "estadisticas"

# Criteria-range headers copied verbatim from the Hospital table (row 1)
$stats.Range("B1").Value = "PACIENTE"
$stats.Range("C1").Value = "FECHA ENTRADA"
$stats.Range("D1").Value = "SEXO"
$stats.Range("E1").Value = "EDAD"
$stats.Range("F1").Value = "VALOR CONSULTA"
$stats.Range("G1").Value = "SINTOMAS"
$stats.Range("H1").Value = "EPS"
$stats.Range("I1").Value = "DIAGNOSTICO"
$stats.Range("J1").Value = "SALA"
$stats.Range("K1").Value = "MEDICO"
$stats.Range("L1").Value = "Nº SALA"

# Criteria row (row 2) - empty cells, but carrying the same number formats as
# the source columns (date for FECHA ENTRADA, currency for VALOR CONSULTA)
$hospital.Range("B2").Copy($stats.Range("C2")) | Out-Null
$stats.Range("C2").ClearContents() | Out-Null
$hospital.Range("E2").Copy($stats.Range("F2")) | Out-Null
$stats.Range("F2").ClearContents() | Out-Null

# Database-function labels and the first worked example (CONTAR / DCOUNT)
$stats.Range("A4").Value = "CONTAR"
$stats.Range("B4").Formula = "=DCOUNT(Hospital!A1:K35,Hospital!E1,B1:L2)"

$stats.Range("A5").Value = "PROMEDIO"
$stats.Range("A6").Value = "MAX"
$stats.Range("A7").Value = "MIN"
$stats.Range("A8").Value = "SUMA"

$stats.Range("A10").Value = "RANGO"
$stats.Range("B10").Value = "A1:K35"
$stats.Range("A11").Value = "Consulta"
$stats.Range("B11").Value = "E1"

# Validation on the criteria row of the new sheet (mirrors Hospital's rules)
$stats.Range("F2").Validation.Add(2, 1, 7, 0) | Out-Null
$stats.Range("E2").Validation.Add(1, 1, 1, 0, 125) | Out-Null

# Column width for the narrow "EPS" criteria column
$stats.Columns.Item(8).AutoFit() | Out-Null

# Selection/active-cell bookkeeping to match the saved state
$hospital.Range("A2:K2").Select() | Out-Null
$stats.Range("F2").Select() | Out-Null

# Make "estadisticas" the active/selected tab, as in the saved workbook
$stats.Select() | Out-Null
